$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1): I1 = "I0", J1 = "IF" ---------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting already used by the other header cells (bold,
# centered horizontally, top-aligned vertically, thin box border).
$hdr = $ws.Range("I1:J1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.LineStyle = 1         # xlContinuous
$hdr.Borders.Weight = 2            # xlThin

# --- New data cells (rows 2-13) for columns I (I0) and J (IF) -----------
$values = @(
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(6, 6),
    @(7, 8),
    @(4, 4),
    @(8, 9),
    @(7, 7),
    @(5, 6),
    @(9, 9),
    @(7, 7),
    @(5, 5)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
